$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead text edits: issue number and week-covering date range ---
$ws.Range("A8").Characters(21, 1).Text = "9"
$ws.Range("C9").Characters(27, 9).Text = "2/27/2023"
$ws.Range("C9").Characters(47, 9).Text = "3/5/2023"

# --- Cells flipping to the "no data" text placeholders (copy value+style from a stable reference cell) ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("C14").Copy($ws.Range("C17"))
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("G28"))
$ws.Range("E14").Copy($ws.Range("H28"))
$ws.Range("C14").Copy($ws.Range("G29"))
$ws.Range("E14").Copy($ws.Range("H29"))
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))

# --- Cells flipping from a text placeholder to a real number (copy numeric style from a stable reference cell, then set value) ---
$ws.Range("K15").Copy($ws.Range("L14"))
$ws.Range("L14").Value = 100
$ws.Range("I14").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 4
$ws.Range("I14").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 3

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -71.428571428571
$ws.Range("H15").Value = -100
$ws.Range("M15").Value = -50
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -41.666666666666
$ws.Range("I16").Value = 21
$ws.Range("K16").Value = -12.5
$ws.Range("L16").Value = 133.333333333333
$ws.Range("M16").Value = -47.5
$ws.Range("N16").Value = -81.578947368421
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -23.076923076923
$ws.Range("I17").Value = 27
$ws.Range("J17").Value = 33
$ws.Range("K17").Value = -18.181818181818
$ws.Range("L17").Value = -15.625
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = -49.056603773584
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 21
$ws.Range("J18").Value = 23
$ws.Range("K18").Value = -8.695652173913
$ws.Range("L18").Value = 10.526315789473
$ws.Range("M18").Value = -44.736842105263
$ws.Range("N18").Value = -88.888888888888
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 266.666666666667
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 54.838709677419
$ws.Range("I19").Value = 101
$ws.Range("J19").Value = 79
$ws.Range("K19").Value = 27.848101265822
$ws.Range("L19").Value = 83.636363636363
$ws.Range("M19").Value = 24.691358024691
$ws.Range("N19").Value = -12.931034482758
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 13
$ws.Range("J20").Value = 23
$ws.Range("K20").Value = -43.478260869565
$ws.Range("L20").Value = -13.333333333333
$ws.Range("M20").Value = -60.606060606060
$ws.Range("N20").Value = -97.560975609756
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 66.666666666666
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = -4.705882352941
$ws.Range("I21").Value = 186
$ws.Range("J21").Value = 184
$ws.Range("K21").Value = 1.086956521739
$ws.Range("L21").Value = 40.909090909090
$ws.Range("M21").Value = -11.004784688995
$ws.Range("N21").Value = -81.728880157170
$ws.Range("G23").Value = 5
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -15.789473684210
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = 14.285714285714
$ws.Range("I24").Value = 199
$ws.Range("J24").Value = 177
$ws.Range("K24").Value = 12.429378531073
$ws.Range("L24").Value = 31.788079470198
$ws.Range("M24").Value = 25.949367088607
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 5
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 12.5
$ws.Range("I25").Value = 48
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = 20
$ws.Range("L25").Value = 20
$ws.Range("M25").Value = -2.040816326530
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -50
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 28.571428571428
$ws.Range("L27").Value = 800
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -20
$ws.Range("N28").Value = -60
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = -25
$ws.Range("N29").Value = -66.666666666666
$ws.Range("G30").Value = 3
$ws.Range("J30").Value = 4
